$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly crypto-list data refresh (prices + 1h volume%) plus a few rows where
# the coin ranking order swapped places (B/C/D/E all change together).
#
# NOTE: column D holds prices as TEXT (e.g. '571.58', '1.00', '62.725.37' using
# '.' as a thousands separator). Excel's COM Value setter auto-coerces plain
# numeric-looking strings into floating point numbers, which would both change
# the cell type and mangle values like '1.00' -> 1 or introduce float noise like
# '571.58' -> 571.58000000000004. Setting NumberFormat to Text ("@") first forces
# the assignment to keep the exact original string.

# Update price (D) and volume (E) columns with the refreshed market data
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.725.37"
$ws.Range("E2").Value = "  +4.98%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.348.03"
$ws.Range("E3").Value = "  +4.82%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "571.58"
$ws.Range("E5").Value = "  +6.93%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.18"
$ws.Range("E6").Value = "  +5.82%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.347.77"
$ws.Range("E8").Value = "  +4.50%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.533"
$ws.Range("E9").Value = "  +1.39%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.44"
$ws.Range("E10").Value = "  +1.78%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.118"
$ws.Range("E11").Value = "  +5.37%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.439"
$ws.Range("E12").Value = "  +3.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.926.93"
$ws.Range("E13").Value = "  +4.79%  "
$ws.Range("E14").Value = "  +0.60%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.91"
$ws.Range("E15").Value = "  +4.71%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000180"
$ws.Range("E16").Value = "  +4.49%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.750.27"
$ws.Range("E17").Value = "  +4.87%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.350.47"
$ws.Range("E18").Value = "  +5.63%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.33"
$ws.Range("E19").Value = "  +2.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.85"
$ws.Range("E20").Value = "  +6.34%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.44"
$ws.Range("E21").Value = "  +2.85%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "384.76"
$ws.Range("E22").Value = "  +4.88%  "
$ws.Range("E23").Value = "  +0.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.534"
$ws.Range("E24").Value = "  +2.64%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "69.97"
$ws.Range("E25").Value = "  +0.46%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.179"
$ws.Range("E26").Value = "  +6.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.26"
$ws.Range("E27").Value = "  +5.44%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0956"
$ws.Range("E28").Value = "  +7.97%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  -0.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.01"
$ws.Range("E30").Value = "  +6.40%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.45"
$ws.Range("E31").Value = "  +6.32%  "
$ws.Range("E34").Value = "  +10.46%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.74"
$ws.Range("E35").Value = "  +2.89%  "
$ws.Range("E36").Value = "  +11.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "159.58"
$ws.Range("E37").Value = "  +1.63%  "
$ws.Range("E38").Value = "  +12.73%  "
$ws.Range("E39").Value = "  +5.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0738"
$ws.Range("E40").Value = "  +5.87%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "40.83"
$ws.Range("E43").Value = "  +3.71%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.04"
$ws.Range("E46").Value = "  +6.11%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "22.11"
$ws.Range("E47").Value = "  +8.63%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.393.99"
$ws.Range("E48").Value = "  +4.86%  "
$ws.Range("E49").Value = "  +0.10%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.32"
$ws.Range("E50").Value = "  +3.38%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "290.43"
$ws.Range("E51").Value = "  +9.65%  "

# Row swaps: the coin ordering changed, so update B/C/D/E for the affected rows
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.60"
$ws.Range("E32").Value = "  +5.78%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.00"
$ws.Range("E33").Value = "  +3.42%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0323"
$ws.Range("E41").Value = "  +10.55%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.786.60"
$ws.Range("E42").Value = "  -0.06%  "
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.28"
$ws.Range("E44").Value = "  +1.52%  "
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.743"
$ws.Range("E45").Value = "  +4.27%  "
